$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.215.34'
$ws.Range("E2").Value = '  +0.25%  '
$ws.Range("D3").Value = '1.859.68'
$ws.Range("E3").Value = '  +0.09%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.26'
$ws.Range("E5").Value = '  +0.78%  '
$ws.Range("E6").Value = '  -0.03%  '
$ws.Range("E7").Value = '  +0.36%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2853'
$ws.Range("E8").Value = '  +1.12%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06518'
$ws.Range("E9").Value = '  -0.77%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.69'
$ws.Range("E10").Value = '  +8.13%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07896'
$ws.Range("E11").Value = '  +0.83%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '97.38'
$ws.Range("E12").Value = '  +0.58%  '
$ws.Range("D13").Value = '1.864.62'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.151'
$ws.Range("E14").Value = '  +0.77%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6787'
$ws.Range("E15").Value = '  +1.91%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '278.98'
$ws.Range("E16").Value = '  -1.19%  '
$ws.Range("D17").Value = '30.211.01'
$ws.Range("E17").Value = '  +0.14%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.58'
$ws.Range("E18").Value = '  +7.62%  '
$ws.Range("E19").Value = '  -0.02%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.370'
$ws.Range("E20").Value = '  -1.30%  '
$ws.Range("D21").Value = '2.108.72'
$ws.Range("E21").Value = '  +0.02%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.000007299'
$ws.Range("E22").Value = '  +0.79%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.163'
$ws.Range("E24").Value = '  +0.44%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '167.24'
$ws.Range("E25").Value = '  -0.43%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.213'
$ws.Range("E26").Value = '  -1.24%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.07'
$ws.Range("E27").Value = '  +0.87%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.929'
$ws.Range("E28").Value = '  +0.66%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.382'
$ws.Range("E29").Value = '  +3.49%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09708'
$ws.Range("E30").Value = '  +1.38%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.362'
$ws.Range("E31").Value = '  -1.32%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.477'
$ws.Range("E32").Value = '  +0.48%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.038'
$ws.Range("E33").Value = '  -1.66%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04723'
$ws.Range("E34").Value = '  +1.15%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.132'
$ws.Range("E35").Value = '  +2.88%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7052'
$ws.Range("E36").Value = '  +0.51%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.707'
$ws.Range("E37").Value = '  +0.31%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01861'
$ws.Range("E38").Value = '  +0.55%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.630'
$ws.Range("E39").Value = '  +4.83%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.329'
$ws.Range("E40").Value = '  +0.05%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '74.44'
$ws.Range("E41").Value = '  +3.20%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.950'
$ws.Range("E42").Value = '  +1.04%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8478'
$ws.Range("E43").Value = '  -0.59%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4169'
$ws.Range("E44").Value = '  +0.34%  '
$ws.Range("E45").Value = '  -0.06%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '103.28'
$ws.Range("E46").Value = '  -0.43%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '973.65'
$ws.Range("E47").Value = '  -1.69%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.171'
$ws.Range("E48").Value = '  -0.87%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.254'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '34.07'
$ws.Range("E51").Value = '  +0.15%  '
